{"js": "// The document stores the element's id as three separate runs:\n//   \"<id>\"  +  \"p052r_1\"  +  \"</id>\"\n// (the middle run carries plain/black formatting while the two tag runs\n// use the Courier New / #7f6000 \"tag\" style). The edit collapses them\n// into a single run \"<id>p052r_1</id>\" that keeps the tag formatting.\n//\n// Locate the run holding the unique id value, then rewrite the text of\n// the whole enclosing paragraph in one shot - Word (and this host) then\n// re-emits a single run carrying the formatting of the first original\n// run, which is exactly the merge the diff shows.\nconst body = context.document.body;\n\nconst hits = body.search(\"p052r_1\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the 'p052r_1' id text in the document\");\n}\n\nconst hitRange = hits.items[0];\nconst paragraph = hitRange.paragraphs.getFirst();\nconst paraRange = paragraph.getRange();\n\n// Rewriting the paragraph's full text in one call forces the host to\n// re-emit it as a single run (carrying the first original run's\n// formatting), collapsing the previous \"<id>\" / \"p052r_1\" / \"</id>\"\n// three-run split into one \"<id>p052r_1</id>\" run - exactly the merge\n// the diff shows.\nparaRange.insertText(\"<id>p052r_1</id>\", \"Replace\");\nawait context.sync();\n", "ps1": "# The element's id is stored as three separate runs:\n#   \"<id>\"  +  \"p052r_1\"  +  \"</id>\"\n# (the middle run is plain/black, the two tag runs use the Courier New /\n# #7f6000 \"tag\" style). The edit collapses them into a single run\n# \"<id>p052r_1</id>\" that keeps the tag formatting.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the id, without assuming a fixed index.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -like \"*p052r_1*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $rng = $target.Range\n\n    # Replacing the whole paragraph range in one Find/Replace pass makes\n    # Word re-emit it as a single run (carrying the first original run's\n    # formatting), merging the previous \"<id>\" / \"p052r_1\" / \"</id>\"\n    # three-run split into one \"<id>p052r_1</id>\" run - exactly the merge\n    # the diff shows.\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute(\"<id>p052r_1</id>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<id>p052r_1</id>\", 2)\n}\n"}
